$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 41, shifting existing rows 41-132 down to 42-133
$ws.Rows(41).Insert()

# Populate the newly inserted row 41 with the new weekly data point
$ws.Range("A41").Value = 8
$ws.Range("B41").Value = "Terminal La Palmera de La Serena"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44519
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 100112037
$ws.Range("G41").Value = "Cebollín"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 900
$ws.Range("L41").Value = 1000
$ws.Range("M41").Value = 950
$ws.Range("N41").Value = "$/paquete 6 unidades"
$ws.Range("O41").Value = "Provincia del Elquí"
$ws.Range("P41").Value = 158
$ws.Range("Q41").Value = 6
$ws.Range("R41").Value = "Hortaliza"
